# Update column G ("K") values on the active sheet to reflect the
# regenerated save_data (K computed instead of Strike#, std/mean
# recalculated, s_vals written out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    16 = 1
    17 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
